$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45937
$ws.Range("B2").Value = 4426.5307593614
$ws.Range("C2").Value = 5450.56345220337
$ws.Range("D2").Value = 7220
$ws.Range("E2").Value = 6495.403429
$ws.Range("F2").Value = 12.4765050767487

$ws.Range("A3").Value = 45938
$ws.Range("B3").Value = 4361.55461745392
$ws.Range("C3").Value = 5233.59498568227
$ws.Range("D3").Value = 3620
$ws.Range("E3").Value = 6041.031066
$ws.Range("F3").Value = 137.211309759514

$ws.Range("A4").Value = 45939
$ws.Range("B4").Value = 4320.84790256862
$ws.Range("C4").Value = 4994.00071930571
$ws.Range("D4").Value = 3620
$ws.Range("E4").Value = 5961.614866
$ws.Range("F4").Value = 125.615320114046

$ws.Range("A5").Value = 45940
$ws.Range("B5").Value = 5326.27534391433
$ws.Range("C5").Value = 4652.20453449099
$ws.Range("D5").Value = 3620
$ws.Range("E5").Value = 7687.700096
$ws.Range("F5").Value = 141.401220274027

$ws.Range("A6").Value = 45941
$ws.Range("B6").Value = 1375.06784482769
$ws.Range("C6").Value = 3101.76005942804
$ws.Range("D6").Value = 3620
$ws.Range("E6").Value = 3268.939353
$ws.Range("F6").Value = 57.3179819833477

$ws.Range("A7").Value = 45942
$ws.Range("B7").Value = 1399.99543622148
$ws.Range("C7").Value = 3204.42682318437
$ws.Range("D7").Value = 3620
$ws.Range("E7").Value = 3764.768721
$ws.Range("F7").Value = 81.2166711651203

$ws.Range("A8").Value = 45943
$ws.Range("B8").Value = 5841.28374289435
$ws.Range("C8").Value = 5611.07623982068
$ws.Range("D8").Value = 3620
$ws.Range("E8").Value = 8951.399833
$ws.Range("F8").Value = 212.549680413597

$ws.Range("A9").Value = 45944
$ws.Range("B9").Value = 5841.28374289435
$ws.Range("C9").Value = 5718.88828711896
$ws.Range("D9").Value = 3620
$ws.Range("E9").Value = 8951.399833
$ws.Range("F9").Value = 217.041849051025

$ws.Range("A10").Value = 45945
$ws.Range("B10").Value = 5841.28374289435
$ws.Range("C10").Value = 5640.56894725707
$ws.Range("D10").Value = 3620
$ws.Range("E10").Value = 8970.950755
$ws.Range("F10").Value = 214.593164973447

$ws.Range("A11").Value = 45946
$ws.Range("B11").Value = 5841.28374289435
$ws.Range("C11").Value = 5687.65574538459
$ws.Range("D11").Value = 3620
$ws.Range("E11").Value = 8970.950755
$ws.Range("F11").Value = 216.555114895427

$ws.Range("A12").Value = 45947
$ws.Range("B12").Value = 5841.28374289435
$ws.Range("C12").Value = 5060.85544858568
$ws.Range("D12").Value = 3620
$ws.Range("E12").Value = 8970.950755
$ws.Range("F12").Value = 190.438435862139

$ws.Range("A13").Value = 45948
$ws.Range("B13").Value = 1742.27770790123
$ws.Range("C13").Value = 3466.31480697635
$ws.Range("D13").Value = 3620
$ws.Range("E13").Value = 4473.79576
$ws.Range("F13").Value = 107.409702461463

$ws.Range("A14").Value = 45949
$ws.Range("B14").Value = 1636.94065696827
$ws.Range("C14").Value = 3392.6573677058
$ws.Range("D14").Value = 3620
$ws.Range("E14").Value = 4358.634659
$ws.Range("F14").Value = 103.931307072397

$ws.Range("A15").Value = 45950
$ws.Range("B15").Value = 6392.95297294923
$ws.Range("C15").Value = 5895.5228997011
$ws.Range("D15").Value = 3620
$ws.Range("E15").Value = 9928.465082
$ws.Range("F15").Value = 242.126458697995

